# Generate Report for Handoff
# Update the status of the "2f7505d0-7278-4f41-9484-17a4c7298420" file from
# "In Translation" to "Ready for handoff" on the Overview, zh-cn, and de-de sheets,
# along with the associated handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (row 3 corresponds to the 2f7505d0... file) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-14-18 10:14:25"

# --- zh-cn sheet (row 3 corresponds to the 2f7505d0... file) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-18 10:14:23"

# --- de-de sheet (row 3 corresponds to the 2f7505d0... file) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-18 10:14:25"
